$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6816443800926208
$ws.Range("B1").Value = 0.9870653748512268
$ws.Range("C1").Value = 1.926884770393372
$ws.Range("D1").Value = 3.284640073776245
$ws.Range("E1").Value = 3.639800071716309
